$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 68, pushing the existing rows 68-104 down to 69-105.
$ws.Rows.Item(68).Insert()

# Populate the new row 68 with the new weekly record.
$ws.Cells.Item(68, 1).Value = 4
$ws.Cells.Item(68, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(68, 3).Value = "Los Lagos"
$ws.Cells.Item(68, 4).Value = 44574
$ws.Cells.Item(68, 5).Value = 10
$ws.Cells.Item(68, 6).Value = 100112022
$ws.Cells.Item(68, 7).Value = "Arveja Verde"
$ws.Cells.Item(68, 8).Value = "Sin especificar"
$ws.Cells.Item(68, 9).Value = "Primera"
$ws.Cells.Item(68, 10).Value = 35
$ws.Cells.Item(68, 11).Value = 30000
$ws.Cells.Item(68, 12).Value = 30000
$ws.Cells.Item(68, 13).Value = 30000
$ws.Cells.Item(68, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(68, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(68, 16).Value = 1200
$ws.Cells.Item(68, 17).Value = 25
$ws.Cells.Item(68, 18).Value = "Hortaliza"

# Match the date number format used by the other rows in column D.
$ws.Cells.Item(68, 4).NumberFormat = $ws.Cells.Item(69, 4).NumberFormat
